$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"0.07140000000000001"
$ws.Range("D2").Value = [double]"1.246"
$ws.Range("C3").Value = [double]"0.0346"
$ws.Range("D3").Value = [double]"1.711"
$ws.Range("C4").Value = [double]"0.0005845"
$ws.Range("D4").Value = [double]"2.071"
$ws.Range("C5").Value = [double]"0.008944000000000001"
$ws.Range("D5").Value = [double]"2.283"
$ws.Range("C6").Value = [double]"0.002097"
$ws.Range("D6").Value = [double]"2.467"
$ws.Range("C7").Value = [double]"4.872e-05"
$ws.Range("D7").Value = [double]"2.63"
$ws.Range("C8").Value = [double]"0.006281"
$ws.Range("D8").Value = [double]"2.677"
$ws.Range("C9").Value = [double]"0.04933"
$ws.Range("D9").Value = [double]"2.71"
$ws.Range("C10").Value = [double]"0.04533"
$ws.Range("D10").Value = [double]"2.742"
$ws.Range("C11").Value = [double]"0.001937"
$ws.Range("D11").Value = [double]"2.781"
$ws.Range("C12").Value = [double]"0.001944"
$ws.Range("D12").Value = [double]"2.792"
$ws.Range("C13").Value = [double]"0.0006522"
$ws.Range("D13").Value = [double]"2.94"
$ws.Range("C14").Value = [double]"0.2507"
$ws.Range("D14").Value = [double]"2.227"
$ws.Range("C15").Value = [double]"0.1076"
$ws.Range("D15").Value = [double]"2.284"
$ws.Range("C16").Value = [double]"4.638e-07"
$ws.Range("D16").Value = [double]"3.349"
$ws.Range("C17").Value = [double]"8.296000000000001e-05"
$ws.Range("D17").Value = [double]"3.802"
$ws.Range("C18").Value = [double]"0.0003701"
$ws.Range("D18").Value = [double]"4.139"
$ws.Range("C19").Value = [double]"0.008718999999999999"
$ws.Range("D19").Value = [double]"3.081"
$ws.Range("C20").Value = [double]"0.0006609"
$ws.Range("D20").Value = [double]"4.45"
$ws.Range("C21").Value = [double]"0.007643"
$ws.Range("D21").Value = [double]"4.945"
$ws.Range("C22").Value = [double]"0.003798"
$ws.Range("D22").Value = [double]"5.237"
$ws.Range("C23").Value = [double]"0.0002045"
$ws.Range("D23").Value = [double]"5.45"
$ws.Range("C24").Value = [double]"0.0006052"
$ws.Range("D24").Value = [double]"5.654"
$ws.Range("C25").Value = [double]"0.02814"
$ws.Range("D25").Value = [double]"4.409"
$ws.Range("C26").Value = [double]"0.04602"
$ws.Range("D26").Value = [double]"4.106"
$ws.Range("C27").Value = [double]"0.00139"
$ws.Range("D27").Value = [double]"5.493"
$ws.Range("C28").Value = [double]"0.001937"
$ws.Range("D28").Value = [double]"5.754"
$ws.Range("C29").Value = [double]"0.0002671"
$ws.Range("D29").Value = [double]"5.768"
$ws.Range("C30").Value = [double]"0.0001061"
$ws.Range("D30").Value = [double]"5.735"
$ws.Range("C31").Value = [double]"0.0002671"
$ws.Range("D31").Value = [double]"5.705"
$ws.Range("C32").Value = [double]"0.0004712"
$ws.Range("D32").Value = [double]"5.723"
$ws.Range("C33").Value = [double]"5.681e-06"
$ws.Range("D33").Value = [double]"5.772"
$ws.Range("C34").Value = [double]"0.01732"
$ws.Range("D34").Value = [double]"4.946"
$ws.Range("C35").Value = [double]"0.002922"
$ws.Range("D35").Value = [double]"5.832"
$ws.Range("C36").Value = [double]"0.0006306"
$ws.Range("D36").Value = [double]"6.133"
$ws.Range("C37").Value = [double]"0.0001719"
$ws.Range("D37").Value = [double]"6.269"
$ws.Range("C38").Value = [double]"0.001069"
$ws.Range("D38").Value = [double]"6.361"
$ws.Range("C39").Value = [double]"0.0002812"
$ws.Range("D39").Value = [double]"6.393"
$ws.Range("C40").Value = [double]"0.0001461"
$ws.Range("D40").Value = [double]"6.374"
$ws.Range("C41").Value = [double]"0.0286"
$ws.Range("D41").Value = [double]"5.206"
$ws.Range("C42").Value = [double]"0.000505"
$ws.Range("D42").Value = [double]"5.988"
$ws.Range("C43").Value = [double]"0.4081"
$ws.Range("D43").Value = [double]"4.755"
$ws.Range("C44").Value = [double]"0.002354"
$ws.Range("D44").Value = [double]"6.054"
$ws.Range("C45").Value = [double]"0.00129"
$ws.Range("D45").Value = [double]"6.508"
$ws.Range("C46").Value = [double]"0.00479"
$ws.Range("D46").Value = [double]"6.827"
$ws.Range("C47").Value = [double]"3.551e-05"
$ws.Range("D47").Value = [double]"7.108"
$ws.Range("C48").Value = [double]"8.609e-05"
$ws.Range("D48").Value = [double]"7.314"
$ws.Range("C49").Value = [double]"0.001841"
$ws.Range("D49").Value = [double]"7.351"
$ws.Range("C50").Value = [double]"0.08749"
$ws.Range("D50").Value = [double]"5.858"
$ws.Range("C51").Value = [double]"0.003351"
$ws.Range("D51").Value = [double]"6.655"
$ws.Range("C52").Value = [double]"0.0003045"
$ws.Range("D52").Value = [double]"6.531"
$ws.Range("C53").Value = [double]"0.001467"
$ws.Range("D53").Value = [double]"6.187"
$ws.Range("C54").Value = [double]"0.001435"
$ws.Range("D54").Value = [double]"5.854"
$ws.Range("C55").Value = [double]"0.003602"
$ws.Range("D55").Value = [double]"5.616"
$ws.Range("C56").Value = [double]"0.002066"
$ws.Range("D56").Value = [double]"5.467"
$ws.Range("C57").Value = [double]"0.001003"
$ws.Range("D57").Value = [double]"5.405"
$ws.Range("C58").Value = [double]"0.008064999999999999"
$ws.Range("D58").Value = [double]"5.353"
$ws.Range("C59").Value = [double]"0.006485"
$ws.Range("D59").Value = [double]"5.253"
$ws.Range("C60").Value = [double]"5.74e-05"
$ws.Range("D60").Value = [double]"5.161"
$ws.Range("C61").Value = [double]"0.194"
$ws.Range("D61").Value = [double]"4.912"
$ws.Range("C62").Value = [double]"0.01743"
$ws.Range("D62").Value = [double]"4.294"
$ws.Range("C63").Value = [double]"0.0003016"
$ws.Range("D63").Value = [double]"5.053"
$ws.Range("C64").Value = [double]"0.001136"
$ws.Range("D64").Value = [double]"5.389"
$ws.Range("C65").Value = [double]"0.004142"
$ws.Range("D65").Value = [double]"5.696"
$ws.Range("C66").Value = [double]"0.0009496"
$ws.Range("D66").Value = [double]"6.076"
$ws.Range("C67").Value = [double]"0.001266"
$ws.Range("D67").Value = [double]"6.466"
$ws.Range("C68").Value = [double]"0.002644"
$ws.Range("D68").Value = [double]"6.801"
$ws.Range("C69").Value = [double]"4.899e-06"
$ws.Range("D69").Value = [double]"7.114"
$ws.Range("C70").Value = [double]"0.0002169"
$ws.Range("D70").Value = [double]"7.399"
$ws.Range("C71").Value = [double]"0.000407"
$ws.Range("D71").Value = [double]"7.587"
$ws.Range("C72").Value = [double]"0.4625"
$ws.Range("D72").Value = [double]"5.055"
